$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 449.8
$ws.Range("I33").Value = 277.55554
$ws.Range("K33").Value = 277.55554
$ws.Range("M33").Value = -48.55554000000001
$ws.Range("H64").Value = 6657.387
$ws.Range("I64").Value = 5923.143
$ws.Range("K64").Value = 5923.143
$ws.Range("M64").Value = -5675.143
$ws.Range("H67").Value = 6657.387
$ws.Range("I67").Value = 5923.143
$ws.Range("K67").Value = 5923.143
$ws.Range("M67").Value = -5065.143
$ws.Range("H74").Value = 7514.1
$ws.Range("I74").Value = 5443.8184
$ws.Range("K74").Value = 5443.8184
$ws.Range("M74").Value = -4507.8184
$ws.Range("H76").Value = 4557
$ws.Range("J76").Value = 3278.8
$ws.Range("L76").Value = 3278.8
$ws.Range("N76").Value = -3908.8
$ws.Range("H77").Value = 7514.1
$ws.Range("I77").Value = 5443.8184
$ws.Range("K77").Value = 27219.092
$ws.Range("M77").Value = -22539.092
$ws.Range("H79").Value = 4557
$ws.Range("J79").Value = 3278.8
$ws.Range("L79").Value = 3278.8
$ws.Range("N79").Value = -5462.8
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 19999.715
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 19999.715
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -22495.715
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 19999.715
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 59999.145
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -72479.145
$ws.Range("H132").Value = 933.5714
$ws.Range("I132").Value = 827.1389
$ws.Range("K132").Value = 2481.4167
$ws.Range("M132").Value = 48.58329999999978
$ws.Range("H138").Value = 2993.2034
$ws.Range("J138").Value = 3010.8447
$ws.Range("L138").Value = 9032.534100000001
$ws.Range("N138").Value = -19312.5341

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17254262
$ws.Range("H61").Value = 4557.6445
$ws.Range("I61").Value = 3704.1
$ws.Range("K61").Value = 3704.1
$ws.Range("M61").Value = -3492.1
$ws.Range("H74").Value = 1851.3529
$ws.Range("I74").Value = 1420.5714
$ws.Range("K74").Value = 1420.5714
$ws.Range("M74").Value = -546.5714
$ws.Range("H77").Value = 1851.3529
$ws.Range("I77").Value = 1420.5714
$ws.Range("K77").Value = 7102.857
$ws.Range("M77").Value = -2734.857
$ws.Range("H132").Value = 3345.6853
$ws.Range("I132").Value = 3033.34
$ws.Range("K132").Value = 9100.02
$ws.Range("M132").Value = -6570.02
$ws.Range("H136").Value = 4557.6445
$ws.Range("I136").Value = 3704.1
$ws.Range("K136").Value = 11112.3
$ws.Range("M136").Value = -8562.299999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 408.2
$ws.Range("I22").Value = 397.75
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 397.75
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -224.75
$ws.Range("N22").Value = -796
$ws.Range("H132").Value = 70073.484
$ws.Range("J132").Value = 70073.484
$ws.Range("L132").Value = 70073.484
$ws.Range("N132").Value = -80193.484
$ws.Range("H134").Value = 3679.75
$ws.Range("I134").Value = 2300.9048
$ws.Range("K134").Value = 6902.714399999999
$ws.Range("M134").Value = -4367.714399999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 533.3333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 533.3333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1599.9999
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4095.9999
$ws.Range("H113").Value = 1877.4546
$ws.Range("I113").Value = 1994.8
$ws.Range("K113").Value = 5984.4
$ws.Range("M113").Value = -3814.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 47728.652
$ws.Range("I70").Value = 105818.63
$ws.Range("J70").Value = 5129.3335
$ws.Range("K70").Value = 105818.63
$ws.Range("L70").Value = 5129.3335
$ws.Range("M70").Value = -105548.63
$ws.Range("N70").Value = -5669.3335
$ws.Range("H73").Value = 47728.652
$ws.Range("I73").Value = 105818.63
$ws.Range("J73").Value = 5129.3335
$ws.Range("K73").Value = 105818.63
$ws.Range("L73").Value = 5129.3335
$ws.Range("M73").Value = -104882.63
$ws.Range("N73").Value = -7001.3335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 804.6429000000001
$ws.Range("I22").Value = 766.6
$ws.Range("J22").Value = 899.75
$ws.Range("K22").Value = 766.6
$ws.Range("L22").Value = 899.75
$ws.Range("M22").Value = -471.6
$ws.Range("N22").Value = -1489.75
$ws.Range("H27").Value = 804.6429000000001
$ws.Range("I27").Value = 766.6
$ws.Range("J27").Value = 899.75
$ws.Range("K27").Value = 766.6
$ws.Range("L27").Value = 899.75
$ws.Range("M27").Value = -659.6
$ws.Range("N27").Value = -1113.75
$ws.Range("H46").Value = 8886.433999999999
$ws.Range("I46").Value = 4384.7144
$ws.Range("J46").Value = 10256.521
$ws.Range("K46").Value = 4384.7144
$ws.Range("L46").Value = 10256.521
$ws.Range("M46").Value = -4196.7144
$ws.Range("N46").Value = -10632.521
$ws.Range("H61").Value = 5432
$ws.Range("I61").Value = 3550.2222
$ws.Range("J61").Value = 6971.636
$ws.Range("K61").Value = 3550.2222
$ws.Range("L61").Value = 6971.636
$ws.Range("M61").Value = -3348.2222
$ws.Range("N61").Value = -7375.636
$ws.Range("H113").Value = 5432
$ws.Range("I113").Value = 3550.2222
$ws.Range("J113").Value = 6971.636
$ws.Range("K113").Value = 3550.2222
$ws.Range("L113").Value = 6971.636
$ws.Range("M113").Value = -1380.2222
$ws.Range("N113").Value = -11311.636
$ws.Range("H132").Value = 5684.564
$ws.Range("I132").Value = 5046.6
$ws.Range("J132").Value = 7811.1113
$ws.Range("K132").Value = 15139.8
$ws.Range("L132").Value = 23433.3339
$ws.Range("M132").Value = -12609.8
$ws.Range("N132").Value = -28493.3339
$ws.Range("H133").Value = 78888.5
$ws.Range("J133").Value = 78888.5
$ws.Range("L133").Value = 78888.5
$ws.Range("N133").Value = -83948.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 29334
$ws.Range("J109").Value = 29166.666
$ws.Range("L109").Value = 29166.666
$ws.Range("N109").Value = -31940.666
